$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (H04, H05, H06, H09, H11) for columns E..I on row 1
$ws.Range("E1").Value = "H04"
$ws.Range("F1").Value = "H05"
$ws.Range("G1").Value = "H06"
$ws.Range("H1").Value = "H09"
$ws.Range("I1").Value = "H11"

# Scores for rows 2-15, columns E..I
$ws.Range("E2").Value = 8.5
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 9
$ws.Range("H2").Value = 9.5
$ws.Range("I2").Value = 10

$ws.Range("E3").Value = 8.5
$ws.Range("F3").Value = 9.5
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

$ws.Range("E4").Value = 9.5
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 9.5
$ws.Range("I4").Value = 0

$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0

$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 11
$ws.Range("I6").Value = 10

$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 11
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 10

$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0

$ws.Range("E11").Value = 9.75
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 0

$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 10.5
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0

$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0

# Update the active selection to match the recorded user interaction
$ws.Range("D19").Select()
